$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates
$ws.Range("G4").Value = 1.91
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = 1.57
$ws.Range("P4").Value = 2.25
$ws.Range("Q4").Value = 2.88
$ws.Range("R4").Value = 1.4
$ws.Range("AC4").Value = 5.5
$ws.Range("AE4").Value = 23
$ws.Range("AF4").Value = 101
$ws.Range("AZ4").Value = 126

# Row 13 updates
$ws.Range("H13").Value = 3.2
$ws.Range("I13").Value = 4.85
$ws.Range("J13").Value = 2.35
$ws.Range("K13").Value = 2.02
$ws.Range("S13").Value = 1.44
$ws.Range("T13").Value = 2.42
$ws.Range("U13").Value = 1.93
$ws.Range("V13").Value = 1.7
$ws.Range("W13").Value = 5.6
$ws.Range("X13").Value = 7.3
$ws.Range("AA13").Value = 16
$ws.Range("AB13").Value = 35
$ws.Range("AD13").Value = 6.3
$ws.Range("AE13").Value = 17.5
$ws.Range("AG13").Value = 11.5
$ws.Range("AH13").Value = 28
$ws.Range("AL13").Value = 60
$ws.Range("AM13").Value = 900
$ws.Range("AO13").Value = 8.75
$ws.Range("AP13").Value = 19.5
$ws.Range("AQ13").Value = 32
$ws.Range("AR13").Value = 70
$ws.Range("AS13").Value = 300
$ws.Range("AU13").Value = 7.5
$ws.Range("AV13").Value = 75
$ws.Range("AY13").Value = 32
$ws.Range("BA13").Value = 200
$ws.Range("BB13").Value = 450
